$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 391
$ws.Cells.Item(4, 6).Value = 1690
$ws.Cells.Item(5, 6).Value = 827
$ws.Cells.Item(6, 6).Value = 717
$ws.Cells.Item(7, 6).Value = 2717
$ws.Cells.Item(9, 6).Value = 2096
$ws.Cells.Item(10, 6).Value = 855
$ws.Cells.Item(11, 6).Value = 2344
$ws.Cells.Item(12, 6).Value = 733
$ws.Cells.Item(13, 6).Value = 6746
$ws.Cells.Item(14, 6).Value = 134
$ws.Cells.Item(15, 6).Value = 538
$ws.Cells.Item(16, 6).Value = 1277
$ws.Cells.Item(17, 6).Value = 1537
$ws.Cells.Item(18, 6).Value = 1350
$ws.Cells.Item(19, 6).Value = 1218
$ws.Cells.Item(20, 6).Value = 104
$ws.Cells.Item(21, 6).Value = 2667
$ws.Cells.Item(22, 6).Value = 2156
$ws.Cells.Item(23, 6).Value = 1113
$ws.Cells.Item(24, 6).Value = 1028
$ws.Cells.Item(25, 6).Value = 796
$ws.Cells.Item(26, 6).Value = 1125
$ws.Cells.Item(27, 6).Value = 264
$ws.Cells.Item(28, 6).Value = 5416
$ws.Cells.Item(29, 6).Value = 295
$ws.Cells.Item(30, 6).Value = 1032
$ws.Cells.Item(31, 6).Value = 1282
$ws.Cells.Item(32, 6).Value = 3784
$ws.Cells.Item(33, 6).Value = 642
$ws.Cells.Item(34, 6).Value = 1711
$ws.Cells.Item(35, 6).Value = 1081
$ws.Cells.Item(36, 6).Value = 64
$ws.Cells.Item(37, 6).Value = 287
$ws.Cells.Item(38, 6).Value = 971
$ws.Cells.Item(39, 6).Value = 1062
$ws.Cells.Item(40, 6).Value = 418
$ws.Cells.Item(42, 6).Value = 49
$ws.Cells.Item(43, 6).Value = 111
$ws.Cells.Item(44, 6).Value = 920
$ws.Cells.Item(46, 6).Value = 517
$ws.Cells.Item(49, 6).Value = 62
$ws.Cells.Item(50, 6).Value = 92

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(8, 6).Value = 497
$ws.Cells.Item(10, 6).Value = 405
$ws.Cells.Item(12, 6).Value = 145
$ws.Cells.Item(20, 6).Value = 610
$ws.Cells.Item(21, 6).Value = 268
$ws.Cells.Item(22, 6).Value = 360
$ws.Cells.Item(24, 6).Value = 176
$ws.Cells.Item(28, 6).Value = 311
$ws.Cells.Item(29, 6).Value = 74
$ws.Cells.Item(33, 6).Value = 49
$ws.Cells.Item(37, 6).Value = 214

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 3303
$ws.Cells.Item(5, 6).Value = 406
$ws.Cells.Item(7, 6).Value = 1479
$ws.Cells.Item(9, 6).Value = 412
$ws.Cells.Item(10, 6).Value = 2851
$ws.Cells.Item(11, 6).Value = 326
$ws.Cells.Item(12, 6).Value = 603
$ws.Cells.Item(13, 6).Value = 711
$ws.Cells.Item(14, 6).Value = 1222

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 406
$ws.Cells.Item(3, 6).Value = 1479
$ws.Cells.Item(5, 6).Value = 391
$ws.Cells.Item(6, 6).Value = 412
$ws.Cells.Item(7, 6).Value = 2851
$ws.Cells.Item(8, 6).Value = 1690
$ws.Cells.Item(9, 6).Value = 827
$ws.Cells.Item(10, 6).Value = 2717
$ws.Cells.Item(11, 6).Value = 326
$ws.Cells.Item(13, 6).Value = 855
$ws.Cells.Item(14, 6).Value = 2344
$ws.Cells.Item(15, 6).Value = 6746
$ws.Cells.Item(16, 6).Value = 134
$ws.Cells.Item(17, 6).Value = 603
$ws.Cells.Item(18, 6).Value = 538
$ws.Cells.Item(19, 6).Value = 1277
$ws.Cells.Item(20, 6).Value = 711
$ws.Cells.Item(21, 6).Value = 1538
$ws.Cells.Item(22, 6).Value = 1350
$ws.Cells.Item(23, 6).Value = 1218
$ws.Cells.Item(24, 6).Value = 1222
$ws.Cells.Item(25, 6).Value = 2667
$ws.Cells.Item(26, 6).Value = 2156
$ws.Cells.Item(27, 6).Value = 1113
$ws.Cells.Item(28, 6).Value = 1028
$ws.Cells.Item(29, 6).Value = 796
$ws.Cells.Item(30, 6).Value = 1125
$ws.Cells.Item(31, 6).Value = 264
$ws.Cells.Item(32, 6).Value = 5416
$ws.Cells.Item(33, 6).Value = 295
$ws.Cells.Item(34, 6).Value = 1032
$ws.Cells.Item(35, 6).Value = 1282
$ws.Cells.Item(36, 6).Value = 3784
$ws.Cells.Item(37, 6).Value = 1711
$ws.Cells.Item(38, 6).Value = 1081
$ws.Cells.Item(39, 6).Value = 74
$ws.Cells.Item(40, 6).Value = 64
$ws.Cells.Item(41, 6).Value = 971
$ws.Cells.Item(42, 6).Value = 1062
$ws.Cells.Item(43, 6).Value = 418
$ws.Cells.Item(45, 6).Value = 49
$ws.Cells.Item(46, 6).Value = 920
$ws.Cells.Item(48, 6).Value = 517
$ws.Cells.Item(50, 6).Value = 214
$ws.Cells.Item(51, 6).Value = 92
